$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("phen_oncox")

# Disease Ontology (row 3) source_version: v2023-12-20 -> v2024-01-31
$ws.Range("E3").Value = "v2024-01-31"
